$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 67, shifting rows 67-108 down to 68-109
$ws.Rows.Item(67).Insert()

# Copy the date cell style from the row above (row 66) into new row 67 col D
$ws.Range("D66").Copy()
$ws.Range("D67").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the unchanged columns (copied from the row that used to be at 67, now at 68)
$ws.Range("A67").Value = 11
$ws.Range("B67").Value = "Vega Monumental Concepción"
$ws.Range("C67").Value = "Bíobío"
$ws.Range("D67").Value = 44673
$ws.Range("E67").Value = 8
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100108
$ws.Range("H67").Value = "Tropicales y subtropicales"
$ws.Range("I67").Value = 100108002
$ws.Range("J67").Value = "Mango"
$ws.Range("K67").Value = "Sin especificar"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 200
$ws.Range("N67").Value = 7000
$ws.Range("O67").Value = 7500
$ws.Range("P67").Value = 7250
$ws.Range("Q67").Value = "$/bandeja 4 kilos"
$ws.Range("R67").Value = "Ecuador"
$ws.Range("S67").Value = 1812
$ws.Range("T67").Value = 4
